# Duty cycle and windings are updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Duty Cycle Min (D15): 0.42 -> 0.313
$ws.Range("D15").Value = 0.313

# Duty Cycle Max (D16): 0.21 -> 0.157
$ws.Range("D16").Value = 0.157

# Multiplacation / windings (D24): 0.75 -> 1
$ws.Range("D24").Value = 1

# Update active cell selection to H6 as seen in the diff
$ws.Range("H6").Select()

$excel.Calculate()
